$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header/label fields with actual teacher details
$ws.Range("A3").Value = "নাম: Dr. Md. Hasanuzzaman (Math)"
$ws.Range("A4").Value = "পদবী: অধ্যাপক"
$ws.Range("F5").Value = "বিভাগ :গণিত"

# Enter quantity for Invigilation (row 26) -> 1 unit, price 2700 computed automatically via formula
$ws.Range("G26").Value = 1

# Amount in words
$ws.Range("A32").Value = "কথায়:দুই হাজার সাতশো টাকা মাত্র।"

# Move active selection to B5 (cosmetic, matches final saved state)
$ws.Range("B5").Select()
